$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParaIndexByText($doc, $searchText) {
    $f = $doc.Content.Find
    $f.ClearFormatting()
    $f.Text = $searchText
    $f.Execute() | Out-Null
    if (-not $f.Found) {
        return -1
    }
    $preceding = $doc.Range(0, $f.Parent.Start)
    return $preceding.Paragraphs.Count + 1
}

# --- Step 1: remove the existing "_GoBack" bookmark. It currently sits at
# the end of the "Still struggling..." paragraph; in the target document it
# is recreated inside the brand-new "5/2/16 ..." paragraph instead. ---
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# --- Step 2: find the empty paragraph right before the "Your challenges
# and successes" heading, and replace it with the two new status-update
# paragraphs. ---
$headingIdx = Get-ParaIndexByText $d "Your challenges and successes"
$emptyPara = $d.Paragraphs.Item($headingIdx - 1)

$newBodyXml = "<w:p $wNs>" + `
    "<w:r><w:t xml:space='preserve'>5/2/16 </w:t></w:r>" + `
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" + `
    "<w:bookmarkEnd w:id='0'/>" + `
    "<w:proofErr w:type='gramStart'/>" + `
    "<w:r><w:t>Just</w:t></w:r>" + `
    "<w:proofErr w:type='gramEnd'/>" + `
    "<w:r><w:t xml:space='preserve'> added a new set of features around Proton flux and am getting much better clustering results.</w:t></w:r>" + `
    "</w:p>" + `
    "<w:p $wNs><w:r><w:t>Will add the charts soon.</w:t></w:r></w:p>"

$emptyPara.Range.InsertXML($newBodyXml) | Out-Null

# --- Step 3: move the lastRenderedPageBreak marker from the "Still
# struggling..." paragraph onto the "Your challenges and successes"
# heading paragraph (re-locate via Find since indices shifted). ---
$headingIdx2 = Get-ParaIndexByText $d "Your challenges and successes"
$headingPara = $d.Paragraphs.Item($headingIdx2)

$headingXml = "<w:p $wNs>" + `
    "<w:pPr><w:pStyle w:val='Heading1'/><w:rPr><w:rFonts w:eastAsia='Times New Roman'/></w:rPr></w:pPr>" + `
    "<w:r><w:rPr><w:rFonts w:eastAsia='Times New Roman'/></w:rPr><w:lastRenderedPageBreak/><w:t>Your challenges and successes</w:t></w:r>" + `
    "</w:p>"

$headingPara.Range.InsertXML($headingXml) | Out-Null

$strugglingIdx = Get-ParaIndexByText $d "Still struggling to find the right combination of features for a more accurate clustering model."
$strugglingPara = $d.Paragraphs.Item($strugglingIdx)
$strugglingPara.Range.Text = "Still struggling to find the right combination of features for a more accurate clustering model."

Write-Host "edit applied"
